# Applies the "Added book interaction and modified headphone and notebook
# interaction" edit to the Assets design-document workbook.
#
# Summary of the change on the "Assets" worksheet (sheet2):
#  - The NotebookInteraction row's Description is updated.
#  - Three new C# Script rows are appended: BookInfoInteraction,
#    HeadphonesInteraction and InteractionHintUI.
#  - The whole data range is then sorted ascending by the "Notes" column,
#    which is what the author did after adding the rows (empty "Notes"
#    values sort last, so "Chandelier Light Material" ends up at the
#    bottom).
#  - The "Assets" sheet becomes the active/selected tab, zoomed to 115%,
#    with A4:C5 selected.

$wb = $excel.ActiveWorkbook
$credits = $wb.Worksheets.Item(1)
$assets = $wb.Worksheets.Item(2)

# --- Add the new interaction script rows ------------------------------------
# (Column order below mirrors the original authoring sequence: all three new
# script names first, then their descriptions, then the NotebookInteraction
# description tweak, then the last new description.)
$assets.Range("A12").Value = "Interaction"
$assets.Range("B12").Value = "BookInfoInteraction"
$assets.Range("C12").Value = "C# Script"
$assets.Range("D12").Value = ".cs"
$assets.Range("E12").Value = "Letisja Muco"
$assets.Range("F12").Value = "Owned"
$assets.Range("H12").Value = "With the help of AI"

$assets.Range("A13").Value = "Interaction"
$assets.Range("B13").Value = "HeadphonesInteraction"
$assets.Range("C13").Value = "C# Script"
$assets.Range("D13").Value = ".cs"
$assets.Range("E13").Value = "Letisja Muco"
$assets.Range("F13").Value = "Owned"
$assets.Range("H13").Value = "With the help of AI"

$assets.Range("A14").Value = "Interaction"
$assets.Range("B14").Value = "InteractionHintUI"
$assets.Range("C14").Value = "C# Script"
$assets.Range("D14").Value = ".cs"
$assets.Range("E14").Value = "Letisja Muco"
$assets.Range("F14").Value = "Owned"
$assets.Range("H14").Value = "With the help of AI"

$assets.Range("G12").Value = "Opens and closes the book UI panel so the player can read the note"
$assets.Range("G13").Value = "Plays and pauses music"

# --- Update the existing NotebookInteraction description -------------------
$assets.Range("G3").Value = "Opens+A8:G10 and closes the notebook UI panel so the player can read the note"

$assets.Range("G14").Value = "Gives hint for each interaction"

# --- Sort A2:H14 by column H (Notes), then column B (Name), ascending ------
$sortObj = $assets.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($assets.Range("H1:H14"))
$sortObj.SortFields.Add($assets.Range("B1:B14"))
$sortObj.SetRange($assets.Range("A1:H14"))
$sortObj.Header = 1
$sortObj.Apply()

# --- Make the Assets sheet the active / selected tab ------------------------
$assets.Activate()
$excel.ActiveWindow.Zoom = 115
$assets.Range("A4:C5").Select()

$wb.Save()
